# Refresh the crypto price/volume table (GitHub Actions data pull).
# Columns B-D are forced to literal text via a leading apostrophe
# (Excel's quote-prefix) so things like trailing zeros and
# multi-dot numbers ("42.697.90") are preserved verbatim, matching
# the existing inline-string cells already in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.697.90"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "'2.240.37"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'114.71"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").Value = "'276.81"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "'46.48"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "'9.03"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").Value = "'15.30"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "'2.577.00"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "'2.249.56"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "'42.885.43"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "'72.20"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").Value = "'3.01"
$ws.Range("E23").Value = "  +5.88%  "
$ws.Range("D24").Value = "'231.47"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "'9.31"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "'12.12"
$ws.Range("E26").Value = "  +7.07%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "'40.35"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'3.27"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("D31").Value = "'173.56"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'0.0892"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'5.58"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "'4.43"
$ws.Range("E35").Value = "  +11.25%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "'0.0373"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").Value = "'4.65"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'70.97"
$ws.Range("E41").Value = "  -6.61%  "
$ws.Range("B42").Value = "'Celestia"
$ws.Range("C42").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'13.22"
$ws.Range("E42").Value = "  -6.33%  "
$ws.Range("B43").Value = "'Algorand"
$ws.Range("C43").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.233"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("E46").Value = "  -6.67%  "
$ws.Range("D47").Value = "'1.27"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "'8.45"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "'0.0990"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'100.79"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").Value = "'0.643"
$ws.Range("E51").Value = "  +8.08%  "
